# Atualizando o arquivo XLSX
# Applies updated Betfair Back/Lay odds to Sheet1 for 2025-12-17 fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Iskenderunspor vs Adana 1954 FK)
$ws.Range("N2").Value = 1.34
$ws.Range("P2").Value = 1.34
$ws.Range("Q2").Value = 1.01

# Row 3
$ws.Range("N3").Value = 1.29
$ws.Range("P3").Value = 1.28

# Row 4 (FK Radnicki 1923 vs Cukaricki)
$ws.Range("F4").Value = 2.14
$ws.Range("G4").Value = 2.76
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 3.9
$ws.Range("J4").Value = 3.1
$ws.Range("L4").Value = 1.27
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 3.75
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 1.95
$ws.Range("Q4").Value = 1.62
$ws.Range("R4").Value = 1.41
$ws.Range("S4").Value = 2.56
$ws.Range("T4").Value = 1.59
$ws.Range("U4").Value = 2.22
$ws.Range("V4").Value = 1.34
$ws.Range("W4").Value = 1.57

# Row 5 (Young Boys vs Grasshoppers Zurich)
$ws.Range("H5").Value = 5.4
$ws.Range("I5").Value = 5.8
$ws.Range("J5").Value = 4.8
$ws.Range("N5").Value = 6.2
$ws.Range("O5").Value = 1.16
$ws.Range("P5").Value = 2.8
$ws.Range("Q5").Value = 1.48
$ws.Range("R5").Value = 1.75
$ws.Range("S5").Value = 2.2
$ws.Range("T5").Value = 1.61
$ws.Range("U5").Value = 2.48
$ws.Range("X5").Value = 29
$ws.Range("Y5").Value = 980
$ws.Range("AA5").Value = 130
$ws.Range("AC5").Value = 12
$ws.Range("AJ5").Value = 1000
$ws.Range("AN5").Value = 5.7

# Row 6 (FC Zurich vs Lugano)
$ws.Range("I6").Value = 2.68
$ws.Range("J6").Value = 3.6
$ws.Range("L6").Value = 1.33
$ws.Range("N6").Value = 4.4
$ws.Range("P6").Value = 2.16
$ws.Range("R6").Value = 1.46
$ws.Range("U6").Value = 2.34
$ws.Range("V6").Value = 1.59
$ws.Range("W6").Value = 1.5

# Row 7 (Luzern vs FC Basel)
$ws.Range("H7").Value = 2.04
$ws.Range("J7").Value = 3.85
$ws.Range("K7").Value = 4.3
$ws.Range("N7").Value = 5.6
$ws.Range("O7").Value = 1.18
$ws.Range("P7").Value = 2.52
$ws.Range("R7").Value = 1.64
$ws.Range("S7").Value = 2.3
$ws.Range("T7").Value = 1.53
$ws.Range("AC7").Value = 10

# Row 8 (Dundee Utd vs Celtic)
$ws.Range("F8").Value = 7.4
$ws.Range("H8").Value = 1.5
$ws.Range("I8").Value = 1.52
$ws.Range("J8").Value = 4.6
$ws.Range("K8").Value = 5.1
$ws.Range("L8").Value = 1.3
$ws.Range("O8").Value = 1.21
$ws.Range("Q8").Value = 1.63
$ws.Range("R8").Value = 1.58
$ws.Range("S8").Value = 2.56
$ws.Range("T8").Value = 1.8
$ws.Range("U8").Value = 2.14
$ws.Range("V8").Value = 2.92
$ws.Range("X8").Value = 23
$ws.Range("Y8").Value = 10.5
$ws.Range("AA8").Value = 14
$ws.Range("AD8").Value = 10.5
$ws.Range("AI8").Value = 32
$ws.Range("AO8").Value = 6.2
